$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Usecase"
$ws.Range("J1").Value = "Expected Output`n"
$ws.Range("K1").Value = "Actual output"

$ws.Range("K2").Value = "Nykaa discount,500 coupon for book,Vocher for groceries"
$ws.Range("J3").Value = "senior citizen discount"
$ws.Range("K3").Value = "senior citizen discount"
$ws.Range("J4").Value = "senior citizen discount"
$ws.Range("K4").Value = "senior citizen discount"
$ws.Range("J5").Value = "100 voucher fastrack coupon  "
$ws.Range("K5").Value = "100 voucher fastrack coupon  "
$ws.Range("K6").Value = "100 vocher fastrack coupon ,vocher for groceries,500 coupon for books "
$ws.Range("K7").Value = "Senior citizen discount applied"
$ws.Range("J8").Value = "senior citizen discount applied"
$ws.Range("K8").Value = "senior citizen discount applied"
$ws.Range("J9").Value = "nykaa discount"
$ws.Range("K9").Value = "nykaa discount"
$ws.Range("J10").Value = "nykaa discount,500 coupon for book"
$ws.Range("K10").Value = "nykaa discount,500 coupon for book"
$ws.Range("J11").Value = "senior citizen discount"
$ws.Range("K11").Value = "senior citizen discount"
$ws.Range("J12").Value = "senior citizen discount"
$ws.Range("K12").Value = "senior citizen discount"
$ws.Range("J13").Value = "100  voucher in fastrack"
$ws.Range("K13").Value = "100  voucher in fastrack"
$ws.Range("J14").Value = "100 voucher in fastrack,500 coupon for book"
$ws.Range("K14").Value = "100 voucher in fastrack,500 coupon for book"
$ws.Range("J15").Value = "senior citizen discount applied"
$ws.Range("K15").Value = "senior citizen discount applied"

$ws.Range("I16").Value = "Hosteller Female working" + [char]0x00A0 + " < 45"
$ws.Range("J16").Value = "Nykaa discount applied"
$ws.Range("K16").Value = "Nykaa discount applied"
